# Update the worksheet date and the twenty-five two-digit x two-digit
# multiplication problems/answers to the new generated set.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-01 Thursday", "2025-05-02 Friday"),
    @("15×60=900",  "26×35=910"),
    @("25×38=950",  "44×91=4004"),
    @("73×45=3285", "17×21=357"),
    @("66×54=3564", "92×97=8924"),
    @("53×17=901",  "55×19=1045"),
    @("49×69=3381", "77×28=2156"),
    @("84×76=6384", "99×70=6930"),
    @("75×14=1050", "78×94=7332"),
    @("78×43=3354", "77×57=4389"),
    @("64×91=5824", "39×42=1638"),
    @("25×67=1675", "95×25=2375"),
    @("33×76=2508", "92×95=8740"),
    @("45×95=4275", "44×25=1100"),
    @("32×28=896",  "90×30=2700"),
    @("90×74=6660", "94×26=2444"),
    @("86×29=2494", "49×17=833"),
    @("52×60=3120", "74×87=6438"),
    @("43×51=2193", "83×50=4150"),
    @("94×64=6016", "16×85=1360"),
    @("82×13=1066", "86×66=5676"),
    @("30×50=1500", "47×88=4136"),
    @("62×39=2418", "97×82=7954"),
    @("58×36=2088", "17×73=1241"),
    @("81×88=7128", "33×78=2574"),
    @("15×93=1395", "45×11=495")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
